# "Drop in RMI script files" — revert this workbook to the older upstream
# version that the RMI scripts expect to find:
#   1. Drop the ad-hoc "Texas Data" worksheet (and its notes in the shared
#      strings table go away with it).
#   2. Restore HPEbP's natural-gas-reforming efficiency formula (B3) to the
#      earlier calculation that still folded waste heat into the energy
#      balance; every later cell on that row (C3:AI3) just chains off B3 /
#      its neighbour, so they recompute on their own.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$texasSheet = $wb.Worksheets.Item("Texas Data")
$texasSheet.Delete()

$hp = $wb.Worksheets.Item("HPEbP")
$hp.Range("B3").Formula = "=118/(162+2+46)"

$wb.Application.CalculateFull()
